$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "REVENDEDORAS"
$ws.Range("B1").Value = "N.Fs"
$ws.Range("C1").Value = "DESTINO "

# Data rows 2-15
$data = @(
    @("BETILDE FERNANDES RABELO", "225597 ", "Ananás - TO"),
    @("CASSIANADASILVA CONCEICAO", "225667 ", "Wanderlândia - TO"),
    @("CEYJANEMACIEL DA SILVA", "225617 ", "São Bento do Tocantins - TO"),
    @("CLAUDILENELEMOS ALENCAR", "225664 ", "Darcinópolis - TO"),
    @("CLAUDINEIA SILVA ARAUJO", "225782 ", "Palmeiras do Tocantins - TO"),
    @("CLEANE FONSECA SILVA", "225669 ", "Palmeiras do Tocantins - TO"),
    @("CLEOMAREUZÉBIO DOS SANTOS", "225678 ", "Araguanã - TO"),
    @("CREUZA PEREIRA BRANDAO", "225631 ", "Nazaré - TO"),
    @("DARKLEY RIBEIRO DE BRITO DIAS", "225767 ", "Darcinópolis - TO"),
    @("DAVILENE OLIVEIRA DA SILVA CHAVES", "225633 ", "Piraquê - TO"),
    @("DAYANY GONCALVES LIMA", "225729 ", "Wanderlândia - TO"),
    @("DINA MARIA PORTILHO", "225688 ", "Angico - TO"),
    @("DAMASCENO", "", "nan"),
    @("DOMINGAS DA CRUZ SILVA", "225656 ", "Luzinópolis - TO")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $ws.Cells.Item($r, 2).NumberFormat = "@"
        $ws.Cells.Item($r, 2).Value = $row[1]
    } else {
        $ws.Cells.Item($r, 2).Value = ""
    }
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Remove old rows 16-20 (dimension shrinks from A1:C20 to A1:C15)
$ws.Rows("16:20").Delete()
